# Auto-generated edit script: apply numeric corrections to Leve profit calc sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1845.5883
$ws.Range("J17").Value = 1845.5883
$ws.Range("L17").Value = 5536.7649
$ws.Range("N17").Value = -5872.7649
$ws.Range("H58").Value = 2307.4167
$ws.Range("J58").Value = 3485.7144
$ws.Range("L58").Value = 10457.1432
$ws.Range("N58").Value = -10757.1432
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("H74").Value = 5652.4
$ws.Range("I74").Value = 5102.1904
$ws.Range("K74").Value = 5102.1904
$ws.Range("M74").Value = -4166.1904
$ws.Range("H76").Value = 166670160
$ws.Range("I76").Value = 250003140
$ws.Range("K76").Value = 250003140
$ws.Range("M76").Value = -250002825
$ws.Range("H77").Value = 5652.4
$ws.Range("I77").Value = 5102.1904
$ws.Range("K77").Value = 25510.952
$ws.Range("M77").Value = -20830.952
$ws.Range("H79").Value = 166670160
$ws.Range("I79").Value = 250003140
$ws.Range("K79").Value = 250003140
$ws.Range("M79").Value = -250002048
$ws.Range("H98").Value = 5629.0713
$ws.Range("I98").Value = 660.8
$ws.Range("K98").Value = 660.8
$ws.Range("M98").Value = 837.2
$ws.Range("H122").Value = 5629.0713
$ws.Range("I122").Value = 660.8
$ws.Range("K122").Value = 1982.4
$ws.Range("M122").Value = 467.6000000000001
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3004.8298
$ws.Range("I32").Value = 2635.3914
$ws.Range("K32").Value = 2635.3914
$ws.Range("M32").Value = -2348.3914
$ws.Range("H45").Value = 2579.1428
$ws.Range("I45").Value = 2014.1428
$ws.Range("K45").Value = 2014.1428
$ws.Range("M45").Value = -1637.1428
$ws.Range("H63").Value = 3997.5
$ws.Range("I63").Value = 3994.6667
$ws.Range("K63").Value = 3994.6667
$ws.Range("M63").Value = -3308.6667
$ws.Range("H66").Value = 3997.5
$ws.Range("I66").Value = 3994.6667
$ws.Range("K66").Value = 19973.3335
$ws.Range("M66").Value = -16541.3335
$ws.Range("H88").Value = 7969.6
$ws.Range("I88").Value = 1738
$ws.Range("J88").Value = 12124
$ws.Range("K88").Value = 1738
$ws.Range("L88").Value = 12124
$ws.Range("M88").Value = -1332
$ws.Range("N88").Value = -12936
$ws.Range("H91").Value = 7969.6
$ws.Range("I91").Value = 1738
$ws.Range("J91").Value = 12124
$ws.Range("K91").Value = 1738
$ws.Range("L91").Value = 12124
$ws.Range("M91").Value = -334
$ws.Range("N91").Value = -14932
$ws.Range("H122").Value = 1731.8948
$ws.Range("I122").Value = 731.2
$ws.Range("K122").Value = 2193.6
$ws.Range("M122").Value = 256.3999999999996
$ws.Range("H132").Value = 3016.3953
$ws.Range("I132").Value = 3022.205
$ws.Range("K132").Value = 9066.615
$ws.Range("M132").Value = -6536.615
$ws.Range("H133").Value = 139994.5
$ws.Range("J133").Value = 139994.5
$ws.Range("L133").Value = 139994.5
$ws.Range("N133").Value = -145054.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 90000
$ws.Range("J59").Value = 90000
$ws.Range("L59").Value = 90000
$ws.Range("H86").Value = 11707.357
$ws.Range("I86").Value = 4649.8887
$ws.Range("J86").Value = 24410.8
$ws.Range("K86").Value = 4649.8887
$ws.Range("L86").Value = 24410.8
$ws.Range("M86").Value = -3526.8887
$ws.Range("N86").Value = -26656.8
$ws.Range("H89").Value = 11707.357
$ws.Range("I89").Value = 4649.8887
$ws.Range("J89").Value = 24410.8
$ws.Range("K89").Value = 23249.4435
$ws.Range("L89").Value = 122054
$ws.Range("M89").Value = -17633.4435
$ws.Range("N89").Value = -133286
$ws.Range("H105").Value = 4296.143
$ws.Range("I105").Value = 3725.3333
$ws.Range("J105").Value = 4724.25
$ws.Range("K105").Value = 3725.3333
$ws.Range("L105").Value = 4724.25
$ws.Range("M105").Value = -1978.3333
$ws.Range("N105").Value = -8218.25
$ws.Range("N59").Value = -91694

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1235.5333
$ws.Range("I16").Value = 1052.6364
$ws.Range("K16").Value = 1052.6364
$ws.Range("M16").Value = -765.6364000000001
$ws.Range("H31").Value = 2294.6365
$ws.Range("I31").Value = 2082.7
$ws.Range("K31").Value = 2082.7
$ws.Range("M31").Value = -1787.7
$ws.Range("H34").Value = 2294.6365
$ws.Range("I34").Value = 2082.7
$ws.Range("K34").Value = 2082.7
$ws.Range("M34").Value = -1880.7
$ws.Range("H41").Value = 43482.668
$ws.Range("J41").Value = 49989.2
$ws.Range("L41").Value = 49989.2
$ws.Range("N41").Value = -50845.2
$ws.Range("H105").Value = 1597.8
$ws.Range("I105").Value = 999.6667
$ws.Range("J105").Value = 2495
$ws.Range("K105").Value = 999.6667
$ws.Range("L105").Value = 2495
$ws.Range("M105").Value = 747.3333
$ws.Range("N105").Value = -5989
$ws.Range("H113").Value = 1235.5333
$ws.Range("I113").Value = 1052.6364
$ws.Range("K113").Value = 1052.6364
$ws.Range("M113").Value = 1117.3636
$ws.Range("H122").Value = 1399.4
$ws.Range("I122").Value = 726.7273
$ws.Range("J122").Value = 3249.25
$ws.Range("K122").Value = 2180.1819
$ws.Range("L122").Value = 9747.75
$ws.Range("M122").Value = 269.8181
$ws.Range("N122").Value = -14647.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 88.89474
$ws.Range("I12").Value = 133.16667
$ws.Range("J12").Value = 68.46154
$ws.Range("K12").Value = 399.50001
$ws.Range("L12").Value = 205.38462
$ws.Range("M12").Value = -226.50001
$ws.Range("N12").Value = -551.38462
$ws.Range("H38").Value = 290.81818
$ws.Range("I38").Value = 200.83333
$ws.Range("J38").Value = 398.8
$ws.Range("K38").Value = 602.49999
$ws.Range("L38").Value = 1196.4
$ws.Range("M38").Value = -255.49999
$ws.Range("N38").Value = -1890.4
$ws.Range("H39").Value = 12499.5
$ws.Range("J39").Value = 19999
$ws.Range("L39").Value = 59997
$ws.Range("N39").Value = -60585
$ws.Range("H120").Value = 34267.23
$ws.Range("I120").Value = 15912.5
$ws.Range("K120").Value = 47737.5
$ws.Range("M120").Value = -42899.5
$ws.Range("H132").Value = 1277.8823
$ws.Range("J132").Value = 1697.1875
$ws.Range("L132").Value = 15274.6875
$ws.Range("N132").Value = -20334.6875
$ws.Range("H137").Value = 24153.34
$ws.Range("J137").Value = 6462
$ws.Range("L137").Value = 19386
$ws.Range("N137").Value = -29586
$ws.Range("H138").Value = 34495550
$ws.Range("I138").Value = 90913816
$ws.Range("K138").Value = 272741448
$ws.Range("M138").Value = -272736308

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1685.6666
$ws.Range("I113").Value = 1041.6666
$ws.Range("J113").Value = 2973.6667
$ws.Range("K113").Value = 1041.6666
$ws.Range("L113").Value = 2973.6667
$ws.Range("M113").Value = 1128.3334
$ws.Range("N113").Value = -7313.6667
$ws.Range("H126").Value = 5513.0527
$ws.Range("I126").Value = 5129.9
$ws.Range("J126").Value = 5938.778
$ws.Range("K126").Value = 15389.7
$ws.Range("L126").Value = 17816.334
$ws.Range("M126").Value = -12919.7
$ws.Range("N126").Value = -22756.334
$ws.Range("H132").Value = 2534.7273
$ws.Range("I132").Value = 1965.1923
$ws.Range("J132").Value = 4650.143
$ws.Range("K132").Value = 5895.5769
$ws.Range("L132").Value = 13950.429
$ws.Range("M132").Value = -3365.5769
$ws.Range("N132").Value = -19010.429
$ws.Range("H133").Value = 99900
$ws.Range("J133").Value = 99900
$ws.Range("L133").Value = 99900
$ws.Range("N133").Value = -110020
$ws.Range("H134").Value = 48497.5
$ws.Range("J134").Value = 48497.5
$ws.Range("L134").Value = 145492.5
$ws.Range("N134").Value = -150562.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 2028
$ws.Range("I93").Value = 1930.9
$ws.Range("K93").Value = 1930.9
$ws.Range("M93").Value = -682.9000000000001
$ws.Range("H132").Value = 3097.9583
$ws.Range("I132").Value = 2440.0715
$ws.Range("J132").Value = 4019
$ws.Range("K132").Value = 7320.2145
$ws.Range("L132").Value = 12057
$ws.Range("M132").Value = -4790.2145
$ws.Range("N132").Value = -17117

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 32000
$ws.Range("J80").Value = 32000
$ws.Range("L80").Value = 32000
$ws.Range("N80").Value = -33996
$ws.Range("H83").Value = 32000
$ws.Range("J83").Value = 32000
$ws.Range("L83").Value = 96000
$ws.Range("N83").Value = -105984
$ws.Range("H126").Value = 2923.25
$ws.Range("I126").Value = 2248.4167
$ws.Range("K126").Value = 6745.250100000001
$ws.Range("M126").Value = -4275.250100000001
$ws.Range("H132").Value = 1571.6316
$ws.Range("I132").Value = 1579.7646
$ws.Range("K132").Value = 4739.293799999999
$ws.Range("M132").Value = -2209.293799999999
$ws.Range("H138").Value = 89164
$ws.Range("J138").Value = 89164
$ws.Range("L138").Value = 89164
$ws.Range("N138").Value = -99444
$ws.Range("H140").Value = 92649.836
$ws.Range("J140").Value = 92649.836
$ws.Range("L140").Value = 92649.836
$ws.Range("N140").Value = -103009.836
